$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05602"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.556"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.020"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8137"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8389"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06953"
$ws.Range("D11").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09404"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001511"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0005971"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006220"
$ws.Range("D17").Style = "Normal"
$ws.Range("B18").Value = "UpBots"
$ws.Range("C18").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.007490"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17UpBotsUBXTBestin24h"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.499"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "18LEOLEO"
$ws.Range("B20").Value = "BTSEToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.092"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19BTSETokenBTSE"
$ws.Range("B21").Value = "BitpandaEcosystemToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3186"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "20BitpandaEcosystemTokenBEST"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1293"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "21ProBitTokenPROB"
$ws.Range("B23").Value = "MCDex"
$ws.Range("C23").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.738"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "22MCDexMCB"
$ws.Range("B24").Value = "CoinExToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04687"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "23CoinExTokenCET"
$ws.Range("B25").Value = "ZBToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1370"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "24ZBTokenZB"
$ws.Range("B26").Value = "BitKan"
$ws.Range("C26").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.001242"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "25BitKanKAN"
$ws.Range("B27").Value = "HotbitToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.004290"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "26HotbitTokenHTB"
$ws.Range("B28").Value = "NitroEx"
$ws.Range("C28").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.00009702"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "27NitroExNTX"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006201"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1054"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002730"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008168"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005297"
$ws.Range("D45").Style = "Normal"
